$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.06448166666667
$ws.Range("H2").Value = 45.193445
$ws.Range("I2").Value = 0.2352496185839757
$ws.Range("J2").Value = 0.2352496185839758
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.221475333333333
$ws.Range("N2").Value = 3.664426
$ws.Range("O2").Value = 0.1659401365378216
$ws.Range("P2").Value = 0.1659401365378216
$ws.Range("Q2").Value = 18.40089276528555
$ws.Range("R2").Value = 165.60803488757
$ws.Range("S2").Value = 0.03903735382829537
$ws.Range("T2").Value = 0.03903735382829538

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.06448166666667
$ws.Range("H3").Value = 45.193445
$ws.Range("I3").Value = 0.2352496185839757
$ws.Range("J3").Value = 0.2352496185839758
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.457871000000001
$ws.Range("N3").Value = 13.373613
$ws.Range("O3").Value = 0.6056116748500271
$ws.Range("P3").Value = 0.6056116748500272
$ws.Range("Q3").Value = 67.15551595186501
$ws.Range("R3").Value = 604.399643566785
$ws.Range("S3").Value = 0.1424699155184716
$ws.Range("T3").Value = 0.1424699155184716

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.06448166666667
$ws.Range("H4").Value = 45.193445
$ws.Range("I4").Value = 0.2352496185839757
$ws.Range("J4").Value = 0.2352496185839758
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.681593333333333
$ws.Range("N4").Value = 5.04478
$ws.Range("O4").Value = 0.2284481886121514
$ws.Range("P4").Value = 0.2284481886121514
$ws.Range("Q4").Value = 25.33233194078889
$ws.Range("R4").Value = 227.9909874671
$ws.Range("S4").Value = 0.05374234923720875
$ws.Range("T4").Value = 0.05374234923720877

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.39986466666667
$ws.Range("H5").Value = 37.199594
$ws.Range("I5").Value = 0.1936384867313999
$ws.Range("J5").Value = 0.1936384867313999
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.221475333333333
$ws.Range("N5").Value = 3.664426
$ws.Range("O5").Value = 0.1659401365378216
$ws.Range("P5").Value = 0.1659401365378216
$ws.Range("Q5").Value = 15.14612882700489
$ws.Range("R5").Value = 136.315159443044
$ws.Range("S5").Value = 0.03213239692718565
$ws.Range("T5").Value = 0.03213239692718566

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.39986466666667
$ws.Range("H6").Value = 37.199594
$ws.Range("I6").Value = 0.1936384867313999
$ws.Range("J6").Value = 0.1936384867313999
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.457871000000001
$ws.Range("N6").Value = 13.373613
$ws.Range("O6").Value = 0.6056116748500271
$ws.Range("P6").Value = 0.6056116748500272
$ws.Range("Q6").Value = 55.27699710145801
$ws.Range("R6").Value = 497.4929739131221
$ws.Range("S6").Value = 0.1172697282648278
$ws.Range("T6").Value = 0.1172697282648279

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.39986466666667
$ws.Range("H7").Value = 37.199594
$ws.Range("I7").Value = 0.1936384867313999
$ws.Range("J7").Value = 0.1936384867313999
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.681593333333333
$ws.Range("N7").Value = 5.04478
$ws.Range("O7").Value = 0.2284481886121514
$ws.Range("P7").Value = 0.2284481886121514
$ws.Range("Q7").Value = 20.85152975770222
$ws.Range("R7").Value = 187.66376781932
$ws.Range("S7").Value = 0.04423636153938641
$ws.Range("T7").Value = 0.04423636153938643

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.26302666666667
$ws.Range("H8").Value = 39.78908
$ws.Range("I8").Value = 0.2071177776734501
$ws.Range("J8").Value = 0.2071177776734502
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.221475333333333
$ws.Range("N8").Value = 3.664426
$ws.Range("O8").Value = 0.1659401365378216
$ws.Range("P8").Value = 0.1659401365378216
$ws.Range("Q8").Value = 16.20045991867556
$ws.Range("R8").Value = 145.80413926808
$ws.Range("S8").Value = 0.03436915230654249
$ws.Range("T8").Value = 0.03436915230654249

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.26302666666667
$ws.Range("H9").Value = 39.78908
$ws.Range("I9").Value = 0.2071177776734501
$ws.Range("J9").Value = 0.2071177776734502
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.457871000000001
$ws.Range("N9").Value = 13.373613
$ws.Range("O9").Value = 0.6056116748500271
$ws.Range("P9").Value = 0.6056116748500272
$ws.Range("Q9").Value = 59.12486194956001
$ws.Range("R9").Value = 532.1237575460401
$ws.Range("S9").Value = 0.1254329442280337
$ws.Range("T9").Value = 0.1254329442280337

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.26302666666667
$ws.Range("H10").Value = 39.78908
$ws.Range("I10").Value = 0.2071177776734501
$ws.Range("J10").Value = 0.2071177776734502
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.681593333333333
$ws.Range("N10").Value = 5.04478
$ws.Range("O10").Value = 0.2284481886121514
$ws.Range("P10").Value = 0.2284481886121514
$ws.Range("Q10").Value = 22.30301722248889
$ws.Range("R10").Value = 200.7271550024
$ws.Range("S10").Value = 0.04731568113887397
$ws.Range("T10").Value = 0.04731568113887398

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.347281
$ws.Range("H11").Value = 55.041843
$ws.Range("I11").Value = 0.2865143954374152
$ws.Range("J11").Value = 0.2865143954374152
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.221475333333333
$ws.Range("N11").Value = 3.664426
$ws.Range("O11").Value = 0.1659401365378216
$ws.Range("P11").Value = 0.1659401365378216
$ws.Range("Q11").Value = 22.41075117523533
$ws.Range("R11").Value = 201.696760577118
$ws.Range("S11").Value = 0.04754423789893607
$ws.Range("T11").Value = 0.04754423789893608

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.347281
$ws.Range("H12").Value = 55.041843
$ws.Range("I12").Value = 0.2865143954374152
$ws.Range("J12").Value = 0.2865143954374152
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.457871000000001
$ws.Range("N12").Value = 13.373613
$ws.Range("O12").Value = 0.6056116748500271
$ws.Range("P12").Value = 0.6056116748500272
$ws.Range("Q12").Value = 81.78981189875101
$ws.Range("R12").Value = 736.1083070887591
$ws.Range("S12").Value = 0.173516462889496
$ws.Range("T12").Value = 0.173516462889496

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.347281
$ws.Range("H13").Value = 55.041843
$ws.Range("I13").Value = 0.2865143954374152
$ws.Range("J13").Value = 0.2865143954374152
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.681593333333333
$ws.Range("N13").Value = 5.04478
$ws.Range("O13").Value = 0.2284481886121514
$ws.Range("P13").Value = 0.2284481886121514
$ws.Range("Q13").Value = 30.85266541439334
$ws.Range("R13").Value = 277.67398872954
$ws.Range("S13").Value = 0.06545369464898314
$ws.Range("T13").Value = 0.06545369464898317

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.961503666666666
$ws.Range("H14").Value = 14.884511
$ws.Range("I14").Value = 0.077479721573759
$ws.Range("J14").Value = 0.07747972157375901
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.221475333333333
$ws.Range("N14").Value = 3.664426
$ws.Range("O14").Value = 0.1659401365378216
$ws.Range("P14").Value = 0.1659401365378216
$ws.Range("Q14").Value = 6.060354345076222
$ws.Range("R14").Value = 54.543189105686
$ws.Range("S14").Value = 0.01285699557686197
$ws.Range("T14").Value = 0.01285699557686197

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.961503666666666
$ws.Range("H15").Value = 14.884511
$ws.Range("I15").Value = 0.077479721573759
$ws.Range("J15").Value = 0.07747972157375901
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.457871000000001
$ws.Range("N15").Value = 13.373613
$ws.Range("O15").Value = 0.6056116748500271
$ws.Range("P15").Value = 0.6056116748500272
$ws.Range("Q15").Value = 22.117743312027
$ws.Range("R15").Value = 199.059689808243
$ws.Range("S15").Value = 0.04692262394919797
$ws.Range("T15").Value = 0.04692262394919798

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.961503666666666
$ws.Range("H16").Value = 14.884511
$ws.Range("I16").Value = 0.077479721573759
$ws.Range("J16").Value = 0.07747972157375901
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.681593333333333
$ws.Range("N16").Value = 5.04478
$ws.Range("O16").Value = 0.2284481886121514
$ws.Range("P16").Value = 0.2284481886121514
$ws.Range("Q16").Value = 8.343231489175556
$ws.Range("R16").Value = 75.08908340258
$ws.Range("S16").Value = 0.01770010204769907
$ws.Range("T16").Value = 0.01770010204769907
